# Updated via Streamlit Approval System
# Sets APPROVAL_1 (AI) / APPROVAL_2 (AJ) decision plus the
# COST_CENTER/LEDGER_NAME/LEDGER_UNDER/TO/BY (AK:AO) placeholder values
# for each pending row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> approval decision, as applied by the approval workflow.
$decisions = [ordered]@{
    2  = "ACCEPTED"
    3  = "ACCEPTED"
    4  = "ACCEPTED"
    5  = "ACCEPTED"
    6  = "HOLD"
    7  = "HOLD"
    8  = "HOLD"
    9  = "ACCEPTED"
    10 = "ACCEPTED"
    11 = "ACCEPTED"
    12 = "ACCEPTED"
    13 = "ACCEPTED"
    14 = "ACCEPTED"
    15 = "HOLD"
    16 = "ACCEPTED"
    17 = "ACCEPTED"
    18 = "REJECTED"
    19 = "ACCEPTED"
    20 = "HOLD"
    21 = "HOLD"
    22 = "HOLD"
    23 = "HOLD"
    24 = "HOLD"
    25 = "HOLD"
}

foreach ($row in $decisions.Keys) {
    $decision = $decisions[$row]

    # APPROVAL_1 / APPROVAL_2
    $ws.Cells.Item($row, 35).Value = $decision   # AI
    $ws.Cells.Item($row, 36).Value = $decision   # AJ

    # COST_CENTER / LEDGER_NAME / LEDGER_UNDER / TO / BY
    # These are stored as the *text* "0" (not the number 0), so force a
    # text number format first to stop the "0" from being coerced to numeric.
    $placeholderCols = 37, 38, 39, 40, 41   # AK, AL, AM, AN, AO
    foreach ($col in $placeholderCols) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = "0"
    }
}
